$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 2 (diff @ -727)
$ws.Range("H2").Value = 582.4
$ws.Range("I2").Value = 170.88889
$ws.Range("J2").Value = 1199.6666
$ws.Range("K2").Value = 170.88889
$ws.Range("L2").Value = 1199.6666
$ws.Range("M2").Value = -57.88889
$ws.Range("N2").Value = -1425.6666

# Row 21 (diff @ -1673)
$ws.Range("H21").Value = 1000
$ws.Range("I21").Value = 1000
$ws.Range("K21").Value = 1000
$ws.Range("M21").Value = -532

# Row 23 (diff @ -1765)
$ws.Range("H23").Value = 1000
$ws.Range("I23").Value = 1000
$ws.Range("K23").Value = 1000
$ws.Range("M23").Value = -766

# Row 29 (diff @ -2047)
$ws.Range("H29").Value = 2371.6
$ws.Range("I29").Value = 322.5
$ws.Range("K29").Value = 967.5
$ws.Range("M29").Value = -686.5

# Row 38 (diff @ -2503)
$ws.Range("H38").Value = 2544.1428
$ws.Range("I38").Value = 1013.1111
$ws.Range("J38").Value = 5300
$ws.Range("K38").Value = 3039.3333
$ws.Range("L38").Value = 15900
$ws.Range("M38").Value = -2667.3333
$ws.Range("N38").Value = -16644

# Row 51 (diff @ -3155)
$ws.Range("H51").Value = 11320.115
$ws.Range("J51").Value = 11173.28
$ws.Range("L51").Value = 11173.28
$ws.Range("N51").Value = -12141.28

# Row 58 (diff @ -3501)
$ws.Range("H58").Value = 2640.6667
$ws.Range("I58").Value = 261.33334
$ws.Range("J58").Value = 7399.3335
$ws.Range("K58").Value = 784.0000200000001
$ws.Range("L58").Value = 22198.0005
$ws.Range("M58").Value = -634.0000200000001
$ws.Range("N58").Value = -22498.0005

# Row 132 (diff @ -7184)
$ws.Range("H132").Value = 2802.2693
$ws.Range("I132").Value = 3034.6843
$ws.Range("J132").Value = 2171.4285
$ws.Range("K132").Value = 9104.052899999999
$ws.Range("L132").Value = 6514.2855
$ws.Range("M132").Value = -6574.052899999999
$ws.Range("N132").Value = -11574.2855

# Row 135 (diff @ -7331)
$ws.Range("H135").Value = 300001280
$ws.Range("I135").Value = 166667310
$ws.Range("K135").Value = 1500005790
$ws.Range("M135").Value = -1500003255

# Row 138 (diff @ -7481)
$ws.Range("H138").Value = 5003.4116
$ws.Range("I138").Value = 1785.6428
$ws.Range("J138").Value = 6220.946
$ws.Range("K138").Value = 5356.928400000001
$ws.Range("L138").Value = 18662.838
$ws.Range("M138").Value = -216.9284000000007
$ws.Range("N138").Value = -28942.838


# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 26 (diff @ -8962)
$ws.Range("H26").Value = 3000
$ws.Range("I26").Value = 3000
$ws.Range("K26").Value = 3000
$ws.Range("M26").Value = -2670

# Row 46 (diff @ -9933)
$ws.Range("H46").Value = 21998.666
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 21998.666
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 21998.666
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -22636.666

# Row 97 (diff @ -12363)
$ws.Range("H97").Value = 752.17645
$ws.Range("I97").Value = 677.7857
$ws.Range("K97").Value = 677.7857
$ws.Range("M97").Value = -181.7857


# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 20 (diff @ -15541)
$ws.Range("H20").Value = 2718.8823
$ws.Range("J20").Value = 3968
$ws.Range("L20").Value = 3968
$ws.Range("N20").Value = -4462

# Row 74 (diff @ -18130)
$ws.Range("H74").Value = 27994
$ws.Range("J74").Value = 27994
$ws.Range("L74").Value = 27994
$ws.Range("N74").Value = -29866

# Row 77 (diff @ -18274)
$ws.Range("H77").Value = 27994
$ws.Range("J77").Value = 27994
$ws.Range("L77").Value = 83982
$ws.Range("N77").Value = -93342

# Row 80 (diff @ -18418)
$ws.Range("H80").Value = 383.53333
$ws.Range("I80").Value = 709.5
$ws.Range("J80").Value = 265
$ws.Range("K80").Value = 709.5
$ws.Range("L80").Value = 265
$ws.Range("M80").Value = 288.5
$ws.Range("N80").Value = -2261

# Row 83 (diff @ -18568)
$ws.Range("H83").Value = 383.53333
$ws.Range("I83").Value = 709.5
$ws.Range("J83").Value = 265
$ws.Range("K83").Value = 3547.5
$ws.Range("L83").Value = 1325
$ws.Range("M83").Value = 1444.5
$ws.Range("N83").Value = -11309

# Row 94 (diff @ -19113)
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("N94").ClearContents()

# Row 134 (diff @ -21043)
$ws.Range("H134").Value = 34002044
$ws.Range("I134").Value = 36430580
$ws.Range("K134").Value = 109291740
$ws.Range("M134").Value = -109289205


# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 7 (diff @ -21789)
$ws.Range("H7").Value = 200.4
$ws.Range("I7").Value = 80.57143000000001
$ws.Range("K7").Value = 80.57143000000001
$ws.Range("M7").Value = 32.42856999999999

# Row 32 (diff @ -23020)
$ws.Range("H32").Value = 4584
$ws.Range("I32").Value = 4417.8
$ws.Range("J32").Value = 4999.5
$ws.Range("K32").Value = 4417.8
$ws.Range("L32").Value = 4999.5
$ws.Range("M32").Value = -4101.8
$ws.Range("N32").Value = -5631.5

# Row 122 (diff @ -27355)
$ws.Range("H122").Value = 2981.2693
$ws.Range("I122").Value = 2840.56
$ws.Range("K122").Value = 8521.68
$ws.Range("M122").Value = -6071.68

# Row 134 (diff @ -27946)
$ws.Range("H134").Value = 11907848
$ws.Range("I134").Value = 13891878
$ws.Range("K134").Value = 41675634
$ws.Range("M134").Value = -41673099


# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 2 (diff @ -28429)
$ws.Range("H2").Value = 74.833336
$ws.Range("I2").Value = 54.909092
$ws.Range("K2").Value = 329.454552
$ws.Range("M2").Value = -216.454552

# Row 11 (diff @ -28891)
$ws.Range("H11").Value = 124398.375
$ws.Range("I11").Value = 135700.6
$ws.Range("J11").Value = 74
$ws.Range("K11").Value = 407101.8
$ws.Range("L11").Value = 222
$ws.Range("M11").Value = -406961.8
$ws.Range("N11").Value = -502

# Row 23 (diff @ -29497)
$ws.Range("H23").Value = 1136
$ws.Range("I23").Value = 508.75
$ws.Range("J23").Value = 1449.625
$ws.Range("K23").Value = 1526.25
$ws.Range("L23").Value = 4348.875
$ws.Range("M23").Value = -1291.25
$ws.Range("N23").Value = -4818.875

# Row 57 (diff @ -31175)
$ws.Range("H57").Value = 11593.75
$ws.Range("J57").Value = 14308.333
$ws.Range("L57").Value = 42924.999
$ws.Range("N57").Value = -44042.999

# Row 58 (diff @ -31227)
$ws.Range("H58").Value = 8373.5
$ws.Range("J58").Value = 10998.333
$ws.Range("L58").Value = 32994.999
$ws.Range("N58").Value = -33250.999


# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 32 (diff @ -37012)
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()

# Row 45 (diff @ -37643)
$ws.Range("H45").Value = 97325.664
$ws.Range("J45").Value = 97325.664
$ws.Range("L45").Value = 97325.664
$ws.Range("N45").Value = -98443.664

# Row 52 (diff @ -37986)
$ws.Range("H52").Value = 29333
$ws.Range("I52").Value = 29333
$ws.Range("K52").Value = 29333
$ws.Range("M52").Value = -29074

# Row 122 (diff @ -41338)
$ws.Range("H122").Value = 43279.434
$ws.Range("I122").Value = 53886.348
$ws.Range("K122").Value = 161659.044
$ws.Range("M122").Value = -159209.044


# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 16 (diff @ -43074)
$ws.Range("H16").Value = 4922.222
$ws.Range("I16").Value = 4500
$ws.Range("K16").Value = 4500
$ws.Range("M16").Value = -4330

# Row 20 (diff @ -43270)
$ws.Range("H20").Value = 50000
$ws.Range("I20").Value = 50000
$ws.Range("K20").Value = 50000
$ws.Range("M20").Value = -49774

# Row 22 (diff @ -43368)
$ws.Range("H22").Value = 2695.8
$ws.Range("I22").Value = 3028.5715
$ws.Range("K22").Value = 3028.5715
$ws.Range("M22").Value = -2733.5715

# Row 27 (diff @ -43610)
$ws.Range("H27").Value = 2695.8
$ws.Range("I27").Value = 3028.5715
$ws.Range("K27").Value = 3028.5715
$ws.Range("M27").Value = -2921.5715

# Row 46 (diff @ -44529)
$ws.Range("H46").Value = 866.6667
$ws.Range("I46").Value = 822
$ws.Range("K46").Value = 822
$ws.Range("M46").Value = -634

# Row 82 (diff @ -46260)
$ws.Range("H82").Value = 1744.7273
$ws.Range("I82").Value = 1866.1111
$ws.Range("J82").Value = 1198.5
$ws.Range("K82").Value = 1866.1111
$ws.Range("L82").Value = 1198.5
$ws.Range("M82").Value = -1505.1111
$ws.Range("N82").Value = -1920.5

# Row 85 (diff @ -46404)
$ws.Range("H85").Value = 1744.7273
$ws.Range("I85").Value = 1866.1111
$ws.Range("J85").Value = 1198.5
$ws.Range("K85").Value = 1866.1111
$ws.Range("L85").Value = 1198.5
$ws.Range("M85").Value = -618.1111000000001
$ws.Range("N85").Value = -3694.5

# Row 94 (diff @ -46830)
$ws.Range("H94").Value = 4000
$ws.Range("J94").Value = 4000
$ws.Range("L94").Value = 4000
$ws.Range("N94").Value = -5352


# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 34 (diff @ -50766)
$ws.Range("H34").Value = 28000
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 28000
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 28000
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -28406

